$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "SUPPORTED BY"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "2023-06-27"
$ws.Range("E2").Value = "Hurted"
$ws.Range("F2").Value = "Santhosh"
